$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (matching the source formatting,
# e.g. trailing zeros like "1.00") by setting those specific cells to Text format
# before assigning their values.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated values (prices, volume %, and reordered coin rows).
$ws.Range("D2").Value = '51.074.42'
$ws.Range("E2").Value = '  -0.06%  '
$ws.Range("D3").Value = '2.941.45'
$ws.Range("E3").Value = '  +1.29%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '377.61'
$ws.Range("E5").Value = '  +1.99%  '
$ws.Range("D6").Value = '104.32'
$ws.Range("E6").Value = '  +1.41%  '
$ws.Range("D7").Value = '0.541'
$ws.Range("E7").Value = '  +0.28%  '
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").Value = '0.590'
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("D10").Value = '36.84'
$ws.Range("E10").Value = '  -0.15%  '
$ws.Range("E11").Value = '  +0.38%  '
$ws.Range("D12").Value = '0.0838'
$ws.Range("E12").Value = '  +0.38%  '
$ws.Range("D13").Value = '3.406.49'
$ws.Range("E13").Value = '  +1.31%  '
$ws.Range("D14").Value = '18.29'
$ws.Range("E14").Value = '  -0.81%  '
$ws.Range("D15").Value = '7.45'
$ws.Range("E15").Value = '  +1.39%  '
$ws.Range("D16").Value = '2.943.96'
$ws.Range("E16").Value = '  +1.61%  '
$ws.Range("D17").Value = '0.955'
$ws.Range("E17").Value = '  +1.65%  '
$ws.Range("D18").Value = '51.082.28'
$ws.Range("E18").Value = '  +0.09%  '
$ws.Range("E19").Value = '  +1.15%  '
$ws.Range("D20").Value = '7.32'
$ws.Range("E20").Value = '  +1.02%  '
$ws.Range("D21").Value = '12.82'
$ws.Range("E21").Value = '  -0.62%  '
$ws.Range("D22").Value = '0.0₃0956'
$ws.Range("E22").Value = '  +1.36%  '
$ws.Range("D23").Value = '68.93'
$ws.Range("E23").Value = '  +0.89%  '
$ws.Range("D24").Value = '260.42'
$ws.Range("D25").Value = '2.80'
$ws.Range("E25").Value = '  +3.74%  '
$ws.Range("B26").Value = 'RenderToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D26").Value = '7.23'
$ws.Range("E26").Value = '  +20.18%  '
$ws.Range("B27").Value = 'Filecoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D27").Value = '7.47'
$ws.Range("E27").Value = '  +5.36%  '
$ws.Range("B28").Value = 'Kaspa'
$ws.Range("C28").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D28").Value = '0.168'
$ws.Range("E28").Value = '  -0.28%  '
$ws.Range("B29").Value = 'Dai'
$ws.Range("C29").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  -0.07%  '
$ws.Range("D30").Value = '25.77'
$ws.Range("E30").Value = '  -0.01%  '
$ws.Range("E31").Value = '  +7.94%  '
$ws.Range("D32").Value = '9.79'
$ws.Range("E32").Value = '  -1.01%  '
$ws.Range("D33").Value = '34.51'
$ws.Range("E33").Value = '  -1.26%  '
$ws.Range("E34").Value = '  -2.49%  '
$ws.Range("D35").Value = '50.79'
$ws.Range("E35").Value = '  -0.78%  '
$ws.Range("D36").Value = '0.0445'
$ws.Range("E36").Value = '  +6.21%  '
$ws.Range("E37").Value = '  +0.05%  '
$ws.Range("E38").Value = '  -0.25%  '
$ws.Range("D39").Value = '17.18'
$ws.Range("E39").Value = '  +1.31%  '
$ws.Range("E40").Value = '  -2.77%  '
$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D41").Value = '1.83'
$ws.Range("E41").Value = '  -0.77%  '
$ws.Range("B42").Value = 'Stellar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D42").Value = '0.115'
$ws.Range("E42").Value = '  +2.12%  '
$ws.Range("D43").Value = '122.34'
$ws.Range("E43").Value = '  +4.10%  '
$ws.Range("D44").Value = '21.92'
$ws.Range("E44").Value = '  -0.72%  '
$ws.Range("B45").Value = 'TheGraph'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D45").Value = '0.280'
$ws.Range("E45").Value = '  +17.74%  '
$ws.Range("B46").Value = 'WEMIXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").Value = '2.07'
$ws.Range("E46").Value = '  -1.08%  '
$ws.Range("D47").Value = '2.36'
$ws.Range("E47").Value = '  +1.80%  '
$ws.Range("D48").Value = '2.030.67'
$ws.Range("E48").Value = '  -0.50%  '
$ws.Range("E49").Value = '  +0.85%  '
$ws.Range("E50").Value = '  +10.26%  '
$ws.Range("D51").Value = '1.28'
$ws.Range("E51").Value = '  +1.06%  '
